$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 21.69
$ws.Range("P2").Value = 47.38
$ws.Range("Q2").Value = 69.07

$ws.Range("O3").Value = 21.3
$ws.Range("P3").Value = 46.14
$ws.Range("Q3").Value = 67.44

$ws.Range("O6").Value = 22.98
$ws.Range("P6").Value = 41.75
$ws.Range("Q6").Value = 64.73

$ws.Range("O7").Value = 22.75
$ws.Range("P7").Value = 40.96
$ws.Range("Q7").Value = 63.71

$ws.Range("O10").Value = 20.64
$ws.Range("P10").Value = 48.78
$ws.Range("Q10").Value = 69.42

$ws.Range("O11").Value = 20.24
$ws.Range("P11").Value = 48.51
$ws.Range("Q11").Value = 68.75

$ws.Range("O14").Value = 22.33
$ws.Range("P14").Value = 43.15
$ws.Range("Q14").Value = 65.48

$ws.Range("O15").Value = 22.71
$ws.Range("P15").Value = 41.51
$ws.Range("Q15").Value = 64.22

$ws.Range("O18").Value = 21.48
$ws.Range("P18").Value = 46.7
$ws.Range("Q18").Value = 68.18

$ws.Range("O19").Value = 22
$ws.Range("P19").Value = 43.45
$ws.Range("Q19").Value = 65.45

$ws.Range("O22").Value = 21.61
$ws.Range("P22").Value = 46.52
$ws.Range("Q22").Value = 68.13

$ws.Range("O23").Value = 21.22
$ws.Range("P23").Value = 47.06
$ws.Range("Q23").Value = 68.28

$ws.Range("O26").Value = 23.12
$ws.Range("P26").Value = 42.69
$ws.Range("Q26").Value = 65.81

$ws.Range("O27").Value = 22.8
$ws.Range("P27").Value = 41.73
$ws.Range("Q27").Value = 64.53

$ws.Range("O30").Value = 22.26
$ws.Range("P30").Value = 42.34
$ws.Range("Q30").Value = 64.6

$ws.Range("O31").Value = 22.61
$ws.Range("P31").Value = 40.27
$ws.Range("Q31").Value = 62.88

$ws.Range("O34").Value = 20.78
$ws.Range("P34").Value = 48.3
$ws.Range("Q34").Value = 69.08

$ws.Range("O35").Value = 21.49
$ws.Range("P35").Value = 46.94
$ws.Range("Q35").Value = 68.43

$ws.Range("O38").Value = 23.16
$ws.Range("P38").Value = 41.53
$ws.Range("Q38").Value = 64.69

$ws.Range("O39").Value = 22.97
$ws.Range("P39").Value = 41.05
$ws.Range("Q39").Value = 64.02

$ws.Range("O42").Value = 21.51
$ws.Range("P42").Value = 44.44
$ws.Range("Q42").Value = 65.95

$ws.Range("O43").Value = 21.67
$ws.Range("P43").Value = 43.12
$ws.Range("Q43").Value = 64.79

$ws.Range("O46").Value = 20.67
$ws.Range("P46").Value = 47.79
$ws.Range("Q46").Value = 68.46

$ws.Range("O47").Value = 21.02
$ws.Range("P47").Value = 46.93
$ws.Range("Q47").Value = 67.95

$ws.Range("O50").Value = 21.67
$ws.Range("P50").Value = 44.47
$ws.Range("Q50").Value = 66.14

$ws.Range("O51").Value = 22.71
$ws.Range("P51").Value = 41.92
$ws.Range("Q51").Value = 64.63

$ws.Range("O54").Value = 20.78
$ws.Range("P54").Value = 46.29
$ws.Range("Q54").Value = 67.07

$ws.Range("O55").Value = 21.47
$ws.Range("P55").Value = 46.03
$ws.Range("Q55").Value = 67.5

$ws.Range("O58").Value = 19.62
$ws.Range("P58").Value = 50.21
$ws.Range("Q58").Value = 69.83

$ws.Range("O59").Value = 21.76
$ws.Range("P59").Value = 45.63
$ws.Range("Q59").Value = 67.39

$ws.Range("O62").Value = 21.98
$ws.Range("P62").Value = 44.52
$ws.Range("Q62").Value = 66.5

$ws.Range("O63").Value = 21.44
$ws.Range("P63").Value = 43.38
$ws.Range("Q63").Value = 64.82

$ws.Range("O66").Value = 22.6
$ws.Range("P66").Value = 39.6
$ws.Range("Q66").Value = 62.2

$ws.Range("O67").Value = 22.52
$ws.Range("P67").Value = 40.29
$ws.Range("Q67").Value = 62.81

$ws.Range("O70").Value = 24.17
$ws.Range("P70").Value = 36.89
$ws.Range("Q70").Value = 61.06

$ws.Range("O71").Value = 24.17
$ws.Range("P71").Value = 36.89
$ws.Range("Q71").Value = 61.06

$ws.Range("O72").Value = 24.17
$ws.Range("P72").Value = 36.89
$ws.Range("Q72").Value = 61.06

$ws.Range("O73").Value = 23.06
$ws.Range("P73").Value = 37.81
$ws.Range("Q73").Value = 60.87

$ws.Range("O76").Value = 21.28
$ws.Range("P76").Value = 44.67
$ws.Range("Q76").Value = 65.95

$ws.Range("O77").Value = 20.83
$ws.Range("P77").Value = 44.48
$ws.Range("Q77").Value = 65.31

$ws.Range("O80").Value = 21.1
$ws.Range("P80").Value = 42.02
$ws.Range("Q80").Value = 63.12

$ws.Range("O81").Value = 20.98
$ws.Range("P81").Value = 41.74
$ws.Range("Q81").Value = 62.72
